# Add a "Price" column (D) to the product table and append two new product
# rows (BN-RED / GH-BROWN), matching the "remove-column"/"sort"/"group"
# example output described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("D1").Value = "Price"

# Existing rows 2-4: fill in the new Price values
$ws.Range("D2").Value = 60
$ws.Range("D3").Value = 200
$ws.Range("D4").Value = 1000

# New row 5: BN-RED / B Necklace / 1000
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "BN-RED"
$ws.Range("C5").Value = "B Necklace"
$ws.Range("D5").Value = 1000

# New row 6: GH-BROWN / G Handbag / 300
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "GH-BROWN"
$ws.Range("C6").Value = "G Handbag"
$ws.Range("D6").Value = 300

# Selection moves to D7 (just past the last data row) as in the saved file
[void]$ws.Range("D7").Select()
